$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 546 (shifts existing rows 546-589 down to 547-590,
# and grows the used range from A1:R589 to A1:R590).
$ws.Rows.Item(546).Insert()

# Populate the newly inserted row 546 with the new price record.
$row = 546
$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 45106
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = 100112043
$ws.Cells.Item($row, 7).Value = "Pepino ensalada"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 105
$ws.Cells.Item($row, 11).Value = 10500
$ws.Cells.Item($row, 12).Value = 11000
$ws.Cells.Item($row, 13).Value = 10762
$ws.Cells.Item($row, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 179
$ws.Cells.Item($row, 17).Value = 60
$ws.Cells.Item($row, 18).Value = "Hortaliza"
